$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 966.3333
$ws.Range("I28").Value = 966.3333
$ws.Range("K28").Value = 966.3333
$ws.Range("M28").Value = -481.3333
$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("K76").Value = 2000
$ws.Range("M76").Value = -1685
$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("K79").Value = 2000
$ws.Range("M79").Value = -908
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = $null
$ws.Range("N95").Value = 0
$ws.Range("H137").Value = 5582.8335
$ws.Range("I137").Value = 5249.5
$ws.Range("K137").Value = 15748.5
$ws.Range("M137").Value = -13198.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = $null
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("H26").Value = 5147.3335
$ws.Range("I26").Value = 4576.8
$ws.Range("J26").Value = 8000
$ws.Range("K26").Value = 4576.8
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = -4246.8
$ws.Range("N26").Value = -8660
$ws.Range("H61").Value = 3366.3333
$ws.Range("I61").Value = 3339.8
$ws.Range("K61").Value = 3339.8
$ws.Range("M61").Value = -3127.8
$ws.Range("H74").Value = 2948.5833
$ws.Range("I74").Value = 2408.7896
$ws.Range("K74").Value = 2408.7896
$ws.Range("M74").Value = -1534.7896
$ws.Range("H77").Value = 2948.5833
$ws.Range("I77").Value = 2408.7896
$ws.Range("K77").Value = 12043.948
$ws.Range("M77").Value = -7675.948
$ws.Range("H94").Value = 6789
$ws.Range("J94").Value = 6789
$ws.Range("L94").Value = 6789
$ws.Range("N94").Value = -8591
$ws.Range("H97").Value = 3989.0908
$ws.Range("I97").Value = 3583
$ws.Range("J97").Value = 4699.75
$ws.Range("K97").Value = 3583
$ws.Range("L97").Value = 4699.75
$ws.Range("M97").Value = -3087
$ws.Range("N97").Value = -5691.75
$ws.Range("H102").Value = 1944.875
$ws.Range("I102").Value = 1859.8334
$ws.Range("K102").Value = 1859.8334
$ws.Range("M102").Value = -237.8334
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = $null
$ws.Range("H136").Value = 3366.3333
$ws.Range("I136").Value = 3339.8
$ws.Range("K136").Value = 10019.4
$ws.Range("M136").Value = -7469.400000000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = $null
$ws.Range("H29").Value = 5052.7144
$ws.Range("I29").Value = 92.5
$ws.Range("J29").Value = 11666.333
$ws.Range("K29").Value = 92.5
$ws.Range("L29").Value = 11666.333
$ws.Range("M29").Value = 196.5
$ws.Range("N29").Value = -12244.333
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = $null
$ws.Range("H54").Value = 17020.75
$ws.Range("I54").Value = 11027.667
$ws.Range("J54").Value = 35000
$ws.Range("K54").Value = 11027.667
$ws.Range("L54").Value = 35000
$ws.Range("M54").Value = -10543.667
$ws.Range("N54").Value = -35968
$ws.Range("H82").Value = 13601
$ws.Range("I82").Value = 13601
$ws.Range("K82").Value = 13601
$ws.Range("M82").Value = -13218
$ws.Range("H85").Value = 13601
$ws.Range("I85").Value = 13601
$ws.Range("K85").Value = 13601
$ws.Range("M85").Value = -12275
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = $null

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5809.2354
$ws.Range("I31").Value = 3511.75
$ws.Range("K31").Value = 3511.75
$ws.Range("M31").Value = -3216.75
$ws.Range("H34").Value = 5809.2354
$ws.Range("I34").Value = 3511.75
$ws.Range("K34").Value = 3511.75
$ws.Range("M34").Value = -3309.75
$ws.Range("H58").Value = 1271.6
$ws.Range("J58").Value = 1768.5
$ws.Range("L58").Value = 1768.5
$ws.Range("N58").Value = -2174.5
$ws.Range("H99").Value = 799
$ws.Range("I99").Value = 799
$ws.Range("K99").Value = 799
$ws.Range("M99").Value = 699
$ws.Range("H126").Value = 799
$ws.Range("I126").Value = 799
$ws.Range("K126").Value = 2397
$ws.Range("M126").Value = 73
$ws.Range("H134").Value = 4314
$ws.Range("I134").Value = 4221.75
$ws.Range("J134").Value = 4498.5
$ws.Range("K134").Value = 12665.25
$ws.Range("L134").Value = 13495.5
$ws.Range("M134").Value = -10130.25
$ws.Range("N134").Value = -18565.5
$ws.Range("H136").Value = 1271.6
$ws.Range("J136").Value = 1768.5
$ws.Range("L136").Value = 5305.5
$ws.Range("N136").Value = -10405.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3750
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 15000
$ws.Range("N3").Value = -15224
$ws.Range("H11").Value = 7143253
$ws.Range("I11").Value = 10000369
$ws.Range("K11").Value = 30001107
$ws.Range("M11").Value = -30000967
$ws.Range("H26").Value = 124.5
$ws.Range("I26").Value = 124.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 373.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = $null
$ws.Range("N26").Value = -85.5
$ws.Range("H114").Value = 980.75
$ws.Range("I114").Value = 651.6667
$ws.Range("K114").Value = 1955.0001
$ws.Range("M114").Value = 1298.9999
$ws.Range("H121").Value = 436
$ws.Range("I121").Value = 436
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 1308
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = $null
$ws.Range("N121").Value = 2
$ws.Range("H140").Value = 746.7778
$ws.Range("I140").Value = 746.7778
$ws.Range("K140").Value = 2240.3334
$ws.Range("M140").Value = 2939.6666

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = $null
$ws.Range("N51").Value = 0
$ws.Range("H97").Value = 1677.5
$ws.Range("I97").Value = 266.5
$ws.Range("J97").Value = 4499.5
$ws.Range("K97").Value = 266.5
$ws.Range("L97").Value = 4499.5
$ws.Range("M97").Value = 229.5
$ws.Range("N97").Value = -5491.5
$ws.Range("H113").Value = 1522.6666
$ws.Range("I113").Value = 1522.6666
$ws.Range("K113").Value = 1522.6666
$ws.Range("M113").Value = 647.3334
$ws.Range("H126").Value = 4150
$ws.Range("I126").Value = 4533.3335
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 13600.0005
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -11130.0005
$ws.Range("N126").Value = -13940

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2300
$ws.Range("I68").Value = 2143
$ws.Range("K68").Value = 2143
$ws.Range("M68").Value = -1394
$ws.Range("H71").Value = 2300
$ws.Range("I71").Value = 2143
$ws.Range("K71").Value = 10715
$ws.Range("M71").Value = -6971
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null
$ws.Range("H93").Value = 5250.75
$ws.Range("I93").Value = 5250.75
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 5250.75
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = $null
$ws.Range("N93").Value = -4002.75

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2182.6667
$ws.Range("J107").Value = 2922.2856
$ws.Range("L107").Value = 8766.856800000001
$ws.Range("N107").Value = -12606.8568
$ws.Range("H122").Value = 2598.5264
$ws.Range("I122").Value = 2375.1538
$ws.Range("K122").Value = 7125.4614
$ws.Range("M122").Value = -4675.4614
$ws.Range("H141").Value = 119999.336
$ws.Range("J141").Value = 119999.336
$ws.Range("L141").Value = 119999.336
$ws.Range("N141").Value = -130359.336

Write-Output "Edit complete"